$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4243.5713
$ws.Range("I19").Value = 10253.333
$ws.Range("J19").Value = 1839.6666
$ws.Range("K19").Value = 10253.333
$ws.Range("L19").Value = 1839.6666
$ws.Range("M19").Value = -10078.333
$ws.Range("N19").Value = -2189.6666
$ws.Range("H38").Value = 484.7143
$ws.Range("I38").Value = 65.5
$ws.Range("J38").Value = 3000
$ws.Range("K38").Value = 196.5
$ws.Range("L38").Value = 9000
$ws.Range("M38").Value = 175.5
$ws.Range("N38").Value = -9744
$ws.Range("H132").Value = 3125.5762
$ws.Range("I132").Value = 866.75
$ws.Range("J132").Value = 12982.272
$ws.Range("K132").Value = 2600.25
$ws.Range("L132").Value = 38946.81600000001
$ws.Range("M132").Value = -70.25
$ws.Range("N132").Value = -44006.81600000001
$ws.Range("H137").Value = 2234689.8
$ws.Range("I137").Value = 3517464.2
$ws.Range("J137").Value = 1364235.5
$ws.Range("K137").Value = 10552392.6
$ws.Range("L137").Value = 4092706.5
$ws.Range("M137").Value = -10549842.6
$ws.Range("N137").Value = -4097806.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 651.6053000000001
$ws.Range("I2").Value = 726.0968
$ws.Range("K2").Value = 726.0968
$ws.Range("M2").Value = -613.0968
$ws.Range("H32").Value = 928761.9399999999
$ws.Range("I32").Value = 949756.5600000001
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 949756.5600000001
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -949469.5600000001
$ws.Range("N32").Value = -5574
$ws.Range("H45").Value = 850.6667
$ws.Range("I45").Value = 775.5
$ws.Range("J45").Value = 1001
$ws.Range("K45").Value = 775.5
$ws.Range("L45").Value = 1001
$ws.Range("M45").Value = -398.5
$ws.Range("N45").Value = -1755
$ws.Range("H61").Value = 402594.03
$ws.Range("I61").Value = 346296.9
$ws.Range("J61").Value = 480337.66
$ws.Range("K61").Value = 346296.9
$ws.Range("L61").Value = 480337.66
$ws.Range("M61").Value = -346084.9
$ws.Range("N61").Value = -480761.66
$ws.Range("H74").Value = 14420158
$ws.Range("I74").Value = 11840130
$ws.Range("J74").Value = 18576868
$ws.Range("K74").Value = 11840130
$ws.Range("L74").Value = 18576868
$ws.Range("M74").Value = -11839256
$ws.Range("N74").Value = -18578616
$ws.Range("H77").Value = 14420158
$ws.Range("I77").Value = 11840130
$ws.Range("J77").Value = 18576868
$ws.Range("K77").Value = 59200650
$ws.Range("L77").Value = 92884340
$ws.Range("M77").Value = -59196282
$ws.Range("N77").Value = -92893076
$ws.Range("H102").Value = 6739.091
$ws.Range("I102").Value = 1391.25
$ws.Range("K102").Value = 1391.25
$ws.Range("M102").Value = 230.75
$ws.Range("H110").Value = 1522.4584
$ws.Range("I110").Value = 1614.2941
$ws.Range("J110").Value = 1299.4286
$ws.Range("K110").Value = 1614.2941
$ws.Range("L110").Value = 1299.4286
$ws.Range("M110").Value = 430.7058999999999
$ws.Range("N110").Value = -5389.4286
$ws.Range("H116").Value = 651.6053000000001
$ws.Range("I116").Value = 726.0968
$ws.Range("K116").Value = 726.0968
$ws.Range("M116").Value = 1567.9032
$ws.Range("H136").Value = 402594.03
$ws.Range("I136").Value = 346296.9
$ws.Range("J136").Value = 480337.66
$ws.Range("K136").Value = 1038890.7
$ws.Range("L136").Value = 1441012.98
$ws.Range("M136").Value = -1036340.7
$ws.Range("N136").Value = -1446112.98

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 651.6053000000001
$ws.Range("I3").Value = 726.0968
$ws.Range("K3").Value = 726.0968
$ws.Range("M3").Value = -612.0968
$ws.Range("H86").Value = 5517.1177
$ws.Range("I86").Value = 6928.727
$ws.Range("J86").Value = 2929.1667
$ws.Range("K86").Value = 6928.727
$ws.Range("L86").Value = 2929.1667
$ws.Range("M86").Value = -5805.727
$ws.Range("N86").Value = -5175.1667
$ws.Range("H89").Value = 5517.1177
$ws.Range("I89").Value = 6928.727
$ws.Range("J89").Value = 2929.1667
$ws.Range("K89").Value = 34643.635
$ws.Range("L89").Value = 14645.8335
$ws.Range("M89").Value = -29027.635
$ws.Range("N89").Value = -25877.8335
$ws.Range("H134").Value = 7029.3794
$ws.Range("I134").Value = 9914
$ws.Range("J134").Value = 2942.8333
$ws.Range("K134").Value = 29742
$ws.Range("L134").Value = 8828.499899999999
$ws.Range("M134").Value = -27207
$ws.Range("N134").Value = -13898.4999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 999.9091
$ws.Range("I16").Value = 930.5
$ws.Range("K16").Value = 930.5
$ws.Range("M16").Value = -643.5
$ws.Range("H31").Value = 1858558.5
$ws.Range("I31").Value = 1056.0714
$ws.Range("J31").Value = 2601559.5
$ws.Range("K31").Value = 1056.0714
$ws.Range("L31").Value = 2601559.5
$ws.Range("M31").Value = -761.0714
$ws.Range("N31").Value = -2602149.5
$ws.Range("H34").Value = 1858558.5
$ws.Range("I34").Value = 1056.0714
$ws.Range("J34").Value = 2601559.5
$ws.Range("K34").Value = 1056.0714
$ws.Range("L34").Value = 2601559.5
$ws.Range("M34").Value = -854.0714
$ws.Range("N34").Value = -2601963.5
$ws.Range("H62").Value = 3972237.2
$ws.Range("J62").Value = 4577.7
$ws.Range("L62").Value = 4577.7
$ws.Range("N62").Value = -5825.7
$ws.Range("H65").Value = 3972237.2
$ws.Range("J65").Value = 4577.7
$ws.Range("L65").Value = 22888.5
$ws.Range("N65").Value = -29128.5
$ws.Range("H105").Value = 1146.6666
$ws.Range("I105").Value = 869.41174
$ws.Range("J105").Value = 1820
$ws.Range("K105").Value = 869.41174
$ws.Range("L105").Value = 1820
$ws.Range("M105").Value = 877.58826
$ws.Range("N105").Value = -5314
$ws.Range("H107").Value = 431
$ws.Range("I107").Value = 276.33334
$ws.Range("J107").Value = 843.44446
$ws.Range("K107").Value = 276.33334
$ws.Range("L107").Value = 843.44446
$ws.Range("M107").Value = 1643.66666
$ws.Range("N107").Value = -4683.44446
$ws.Range("H113").Value = 999.9091
$ws.Range("I113").Value = 930.5
$ws.Range("K113").Value = 930.5
$ws.Range("M113").Value = 1239.5
$ws.Range("H134").Value = 17242952
$ws.Range("I134").Value = 26316702
$ws.Range("J134").Value = 2825.7
$ws.Range("K134").Value = 78950106
$ws.Range("L134").Value = 8477.099999999999
$ws.Range("M134").Value = -78947571
$ws.Range("N134").Value = -13547.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 735.93335
$ws.Range("I5").Value = 416.89285
$ws.Range("J5").Value = 1261.4117
$ws.Range("K5").Value = 1250.67855
$ws.Range("L5").Value = 3784.2351
$ws.Range("M5").Value = -1138.67855
$ws.Range("N5").Value = -4008.2351
$ws.Range("H12").Value = 37.714287
$ws.Range("I12").Value = 27.777779
$ws.Range("K12").Value = 83.333337
$ws.Range("M12").Value = 89.666663
$ws.Range("H33").Value = 356.09525
$ws.Range("I33").Value = 270.5
$ws.Range("J33").Value = 470.22223
$ws.Range("K33").Value = 1623
$ws.Range("L33").Value = 2821.33338
$ws.Range("M33").Value = -1340
$ws.Range("N33").Value = -3387.33338
$ws.Range("H68").Value = 890.2030999999999
$ws.Range("I68").Value = 572.8182
$ws.Range("J68").Value = 1056.4524
$ws.Range("K68").Value = 1718.4546
$ws.Range("L68").Value = 3169.357199999999
$ws.Range("M68").Value = -907.4546
$ws.Range("N68").Value = -4791.357199999999
$ws.Range("H71").Value = 890.2030999999999
$ws.Range("I71").Value = 572.8182
$ws.Range("J71").Value = 1056.4524
$ws.Range("K71").Value = 5155.3638
$ws.Range("L71").Value = 9508.071599999999
$ws.Range("M71").Value = -1099.3638
$ws.Range("N71").Value = -17620.0716
$ws.Range("H135").Value = 735.93335
$ws.Range("I135").Value = 416.89285
$ws.Range("J135").Value = 1261.4117
$ws.Range("K135").Value = 3752.03565
$ws.Range("L135").Value = 11352.7053
$ws.Range("M135").Value = -1217.03565
$ws.Range("N135").Value = -16422.7053

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1955940.8
$ws.Range("I132").Value = 3127491.5
$ws.Range("J132").Value = 3356.25
$ws.Range("K132").Value = 9382474.5
$ws.Range("L132").Value = 10068.75
$ws.Range("M132").Value = -9379944.5
$ws.Range("N132").Value = -15128.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 32000
$ws.Range("J123").Value = 32000
$ws.Range("L123").Value = 32000
$ws.Range("N123").Value = -41800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 9163.333000000001
$ws.Range("I51").Value = 980
$ws.Range("J51").Value = 10800
$ws.Range("K51").Value = 980
$ws.Range("L51").Value = 10800
$ws.Range("M51").Value = -470
$ws.Range("N51").Value = -11820
$ws.Range("H107").Value = 689.86664
$ws.Range("I107").Value = 550.8889
$ws.Range("J107").Value = 898.3333
$ws.Range("K107").Value = 1652.6667
$ws.Range("L107").Value = 2694.9999
$ws.Range("M107").Value = 267.3332999999998
$ws.Range("N107").Value = -6534.9999
$ws.Range("H126").Value = 1560.7407
$ws.Range("I126").Value = 1189.3077
$ws.Range("K126").Value = 3567.9231
$ws.Range("M126").Value = -1097.9231
$ws.Range("H132").Value = 20835420
$ws.Range("I132").Value = 33334944
$ws.Range("J132").Value = 2882.4443
$ws.Range("K132").Value = 100004832
$ws.Range("L132").Value = 8647.332900000001
$ws.Range("M132").Value = -100002302
$ws.Range("N132").Value = -13707.3329
$ws.Range("H136").Value = 14766009
$ws.Range("I136").Value = 16576545
$ws.Range("J136").Value = 7403160.5
$ws.Range("K136").Value = 49729635
$ws.Range("L136").Value = 22209481.5
$ws.Range("M136").Value = -49727085
$ws.Range("N136").Value = -22214581.5
